$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.196729183197021
$ws.Range("B1").Value = 3.059962272644043
$ws.Range("C1").Value = 2.623027801513672
$ws.Range("D1").Value = 3.226503610610962
$ws.Range("E1").Value = 3.700998783111572
